# Apply the publications-table update described in the commit:
# add 9 new publication rows (100-108) to Sheet1, and turn the
# citation URL in D104 into a clickable hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100
$ws.Cells.Item(100, 1).Value = "Effect of alcohol on the interleukin 6-mediated inflammatory response in a new mouse model of acute-on-chronic liver injury."
$ws.Cells.Item(100, 2).Value = "Biochim Biophys Acta Mol Basis Dis"
$ws.Cells.Item(100, 4).Value = "Biochim Biophys Acta Mol Basis Dis. 2018 Nov 15;1865(2):298-307. doi: 10.1016/j.bbadis.2018.11.008."
$ws.Cells.Item(100, 5).Value = "['Ersin Karatayli', 'R. A. Hall', 'Susanne Weber', 'Steven Dooley', 'Frank Lammert']"
$ws.Cells.Item(100, 6).Value = "{'data': [{'id': '5', 'type': 'projects'}, {'id': '13', 'type': 'projects'}]}"

# Row 101
$ws.Cells.Item(101, 1).Value = "Liver-specific Repin1 deficiency impairs transient hepatic steatosis in liver regeneration."
$ws.Cells.Item(101, 2).Value = "Sci Rep"
$ws.Cells.Item(101, 3).NumberFormat = "@"
$ws.Cells.Item(101, 3).Value = "2018-01-18"
$ws.Cells.Item(101, 4).Value = "Sci Rep. 2018 Nov 15;8(1):16858. doi: 10.1038/s41598-018-35325-3."
$ws.Cells.Item(101, 5).Value = "['K. Abshagen', 'B. Degenhardt', 'M. Liebig', 'A. Wendt', 'B. Genz', 'U. Schaeper', 'M. Stumvoll', 'Ute Hofmann', 'M. Frank', 'B. Vollmar', 'N. Kloting']"
$ws.Cells.Item(101, 6).Value = "{'data': [{'id': '3', 'type': 'projects'}]}"

# Row 102
$ws.Cells.Item(102, 1).Value = "ABCB4 Gene Aberrations in Human Liver Disease: An Evolving Spectrum"
$ws.Cells.Item(102, 2).Value = "Semin Liver Dis"
$ws.Cells.Item(102, 3).NumberFormat = "@"
$ws.Cells.Item(102, 3).Value = "2018-10-24"
$ws.Cells.Item(102, 4).Value = "Semin Liver Dis 38(04) : 299"
$ws.Cells.Item(102, 5).Value = "['Matthias Reichert', 'Frank Lammert']"
$ws.Cells.Item(102, 6).Value = "{'data': [{'id': '5', 'type': 'projects'}, {'id': '13', 'type': 'projects'}]}"

# Row 103
$ws.Cells.Item(103, 1).Value = "Effects of Gene Variants Controlling Vitamin D Metabolism and Serum Levels on Hepatic Steatosis."
$ws.Cells.Item(103, 2).Value = "Digestion"
$ws.Cells.Item(103, 3).NumberFormat = "@"
$ws.Cells.Item(103, 3).Value = "2018-03-08"
$ws.Cells.Item(103, 4).Value = "Digestion. 2018;97(4):298-308. doi: 10.1159/000485180. Epub 2018 Mar 7."
$ws.Cells.Item(103, 5).Value = "['M. Jamka', 'A. Arslanow', 'A. Bohner', 'M. Krawczyk', 'Susanne Weber', 'F. Grunhage', 'Frank Lammert', 'C. S. Stokes']"
$ws.Cells.Item(103, 6).Value = "{'data': [{'id': '5', 'type': 'projects'}, {'id': '13', 'type': 'projects'}]}"

# Row 104
$ws.Cells.Item(104, 1).Value = "Analytical challenges in human plasma lipidomics: A winding path towards the truth"
$ws.Cells.Item(104, 2).Value = "TrAC Trends in Analytical Chemistry"
$ws.Cells.Item(104, 3).NumberFormat = "@"
$ws.Cells.Item(104, 3).Value = "2018-10-20"
$ws.Cells.Item(104, 5).Value = "['Olga Vvedenskaya', 'Yuting Wang', 'Jacobo Miranda Ackerman', 'Oskar Knittelfelder', 'Andrej Shevchenko']"
$ws.Cells.Item(104, 6).Value = "{'data': [{'id': '3', 'type': 'projects'}]}"

# Row 105
$ws.Cells.Item(105, 1).Value = "Epigenomic map of human liver reveals principles of zonated morphogenic and metabolic control"
$ws.Cells.Item(105, 2).Value = "Nat Commun"
$ws.Cells.Item(105, 3).NumberFormat = "@"
$ws.Cells.Item(105, 3).Value = "2018-12-01"
$ws.Cells.Item(105, 4).Value = "Nat Commun 9(1) : 416"
$ws.Cells.Item(105, 5).Value = "['Mario Brosch', 'Kathrin Kattler', 'Alexander Herrmann', 'Witigo von Schönfels', 'Karl Nordström', 'Daniel Seehofer', 'Georg Damm', 'Thomas Becker', 'Sebastian Zeissig', 'Sophie Nehring', 'Fabian Reichel', 'Vincent Moser', 'Raghavan Veera Thangapandi', 'Felix Stickel', 'Gustavo Baretton', 'Christoph Röcken', 'Michael Muders', 'Madlen Matz-Soja', 'Michael Krawczak', 'Gilles Gasparoni', 'Hella Hartmann', 'Andreas Dahl', 'Clemens Schafmayer', 'Jörn Walter', 'Jochen Hampe']"
$ws.Cells.Item(105, 6).Value = "{'data': [{'id': '9', 'type': 'projects'}]}"

# Row 106
$ws.Cells.Item(106, 1).Value = "The Diurnal Timing of Starvation Differently Impacts Murine Hepatic Gene Expression and Lipid Metabolism – A Systems Biology Analysis Using Self-Organizing Maps"
$ws.Cells.Item(106, 2).Value = "Front. Physiol."
$ws.Cells.Item(106, 3).NumberFormat = "@"
$ws.Cells.Item(106, 3).Value = "2018-09-10"
$ws.Cells.Item(106, 4).Value = "Front. Physiol. 9 : 660"
$ws.Cells.Item(106, 5).Value = "['Christiane Rennert', 'Sebastian Vlaic', 'Eugenia Marbach-Breitrück', 'Carlo Thiel', 'Susanne Sales', 'Andrej Shevchenko', 'Rolf Gebhardt', 'Madlen Matz-Soja']"
$ws.Cells.Item(106, 6).Value = "{'data': [{'id': '9', 'type': 'projects'}]}"

# Row 107
$ws.Cells.Item(107, 1).Value = "Clinical and Functional Relevance of the Monocarboxylate Transporter Family in Disease Pathophysiology and Drug Therapy."
$ws.Cells.Item(107, 2).Value = "Clin Transl Sci"
$ws.Cells.Item(107, 3).NumberFormat = "@"
$ws.Cells.Item(107, 3).Value = "2018-04-17"
$ws.Cells.Item(107, 4).Value = "Clin Transl Sci. 2018 Jul;11(4):352-364. doi: 10.1111/cts.12551. Epub 2018 Apr 16."
$ws.Cells.Item(107, 5).Value = "['P. Fisel', 'E. Schaeffeler', 'Matthias Schwab']"
$ws.Cells.Item(107, 6).Value = "{'data': [{'id': '3', 'type': 'projects'}]}"

# Row 108
$ws.Cells.Item(108, 1).Value = "The fruit fly Drosophila melanogaster as an innovative preclinical ADME model for solute carrier membrane transporters, with consequences for pharmacology and drug therapy."
$ws.Cells.Item(108, 2).Value = "Drug Discov Today"
$ws.Cells.Item(108, 3).NumberFormat = "@"
$ws.Cells.Item(108, 3).Value = "2018-06-12"
$ws.Cells.Item(108, 4).Value = "Drug Discov Today. 2018 Oct;23(10):1746-1760. doi: 10.1016/j.drudis.2018.06.002.  Epub 2018 Jun 8."
$ws.Cells.Item(108, 5).Value = "['Y. Wang', 'B. Moussian', 'E. Schaeffeler', 'Matthias Schwab', 'A. T. Nies']"
$ws.Cells.Item(108, 6).Value = "{'data': [{'id': '3', 'type': 'projects'}]}"

# D104 citation is a URL -- add it as a real hyperlink (this also sets the
# cell text to the URL and applies the built-in Hyperlink style).
$ws.Hyperlinks.Add($ws.Range("D104"), "https://www.sciencedirect.com/science/article/pii/S016599361830428X?via%3Dihub")

